# Generate Report for Archive
#
# The localization-status report is regenerated: the status of the two
# handed-off files moves from "Ready for handoff" to "In Translation" on
# every sheet that tracks it (Overview + each locale sheet), and the
# status/locale columns that previously had to be wide enough to fit the
# longer "Ready for handoff" label are narrowed to fit the new text.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: status is mirrored into the zh-cn / de-de columns (E, F) ---
$overview = $wb.Worksheets.Item("Overview")
[void]$overview.Cells.Replace("Ready for handoff", "In Translation")
$overview.Range("E1").EntireColumn.ColumnWidth = 12.5
$overview.Range("F1").EntireColumn.ColumnWidth = 12.5

# --- Per-locale sheets: status lives in column C ("Status") ---
$zhcn = $wb.Worksheets.Item("zh-cn")
[void]$zhcn.Cells.Replace("Ready for handoff", "In Translation")
$zhcn.Range("C1").EntireColumn.ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
[void]$dede.Cells.Replace("Ready for handoff", "In Translation")
$dede.Range("C1").EntireColumn.ColumnWidth = 12.5
